$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row before row 302 (weekly price update) -- this
# shifts the existing rows 302..371 down to 303..372 and extends the
# sheet's used range to A1:R372.
$ws.Rows.Item(302).Insert()

# Populate the newly inserted row 302 with the new observation.
$ws.Cells.Item(302, 1).Value = 10
$ws.Cells.Item(302, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(302, 3).Value = "La Araucanía"
$ws.Cells.Item(302, 4).Value = 44642
$ws.Cells.Item(302, 5).Value = 9
$ws.Cells.Item(302, 6).Value = 100114014
$ws.Cells.Item(302, 7).Value = "Betarraga"
$ws.Cells.Item(302, 8).Value = "Sin especificar"
$ws.Cells.Item(302, 9).Value = "Primera"
$ws.Cells.Item(302, 10).Value = 45
$ws.Cells.Item(302, 11).Value = 8000
$ws.Cells.Item(302, 12).Value = 8000
$ws.Cells.Item(302, 13).Value = 8000
$ws.Cells.Item(302, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(302, 15).Value = "Región del Maule"
$ws.Cells.Item(302, 16).Value = 667
$ws.Cells.Item(302, 17).Value = 12
$ws.Cells.Item(302, 18).Value = "Hortaliza"
